$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: date changes from 2021-08-30 (44438) to 2021-10-05 (44474); quality Primera -> Especial; volume 100 -> 150
$ws.Range("D27").Value = 44474
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 150

# --- Row 28: date changes from 2021-09-22 (44461) to 2021-08-30 (44438); quality Especial -> Primera; volume 150 -> 100
$ws.Range("D28").Value = 44438
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 100

# --- Row 29: a new weekly record is inserted here (quality Especial, volume 150, prices 30000)
#     Date (D29) stays as it was: 2021-09-22 (44461)
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 150
$ws.Range("N29").Value = 30000
$ws.Range("O29").Value = 30000
$ws.Range("P29").Value = 30000
$ws.Range("S29").Value = 3000

# --- Row 30 (new row): holds what used to be in row 29 before the insertion
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44461
$ws.Range("D30").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107002
$ws.Range("J30").Value = "Chirimoya"
$ws.Range("K30").Value = "Cultivar IV Región"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 25000
$ws.Range("O30").Value = 25000
$ws.Range("P30").Value = 25000
$ws.Range("Q30").Value = "$/bandeja 10 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 2500
$ws.Range("T30").Value = 10
